$d = $word.ActiveDocument

function Replace-ParagraphByAnchor {
    param(
        [string]$AnchorText,
        [string]$NewXml
    )
    $paras = $d.Paragraphs
    $target = $null
    for ($i = 1; $i -le $paras.Count; $i++) {
        $p = $paras.Item($i)
        if ($p.Range.Text.Contains($AnchorText)) {
            $target = $p.Range
            break
        }
    }
    if ($null -eq $target) {
        throw "Anchor not found: $AnchorText"
    }
    $target.InsertXML($NewXml)
}

# 1) View section paragraph: replace the "<insert stuff here>" placeholder
$viewXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:r><w:t xml:space="preserve">The view component was designed to use </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>javafx</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. We created a side panel on the left for the user to control the environment by setting the seed, adding Critters, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>steping</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> through the world, animating the environment, and viewing stats about each type of Critter. As a visual representation of the world the right side of the component is a grid which is a 2D representation of the torus the critters reside on. Each Critter can define its </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>CritterShape</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and the fill and outline colors it would like to be displayed as. We use </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>javafx’s</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> native graphics library to display these shapes and colors in the grid for a presentation that is useful and pleasing.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
Replace-ParagraphByAnchor "insert stuff here" $viewXml

# 2) Controller paragraph: simplify grammar proof-marks, merge runs
$controllerXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:r><w:t>The controls really just involved repackaging most of the commands from Project 4</w:t></w:r><w:r><w:t xml:space="preserve">, except we also had to update the view after many of the “commands” (e.g. spawning new Critters). The main new functionality was animation, which relies on </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>javafx’s</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Timeline class, which is very similar to </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>java.util.Timer</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. Unfortunately, Timer doesn’t work with </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>javafx</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> UI components since the spawned thread doesn’t have access to them, which is why we needed to use Timeline.</w:t></w:r></w:p>'
Replace-ParagraphByAnchor "The controls really just involved" $controllerXml

# 3) Notes bullet about IntelliJ/.class files: fix "si" -> "is", drop stray bookmark
$notes2Xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="NoSpacing"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Because we were developing in IntelliJ instead of Eclipse/Linux, the reflection code that loads the possible Critter classes looks for .java files in the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>src</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> folder instead of .class files in the out folder. If we were to distribute this, we would fix this to look for the .class file</w:t></w:r><w:r><w:t>s instead, which is</w:t></w:r><w:r><w:t xml:space="preserve"> more “correct.” To this end, we made our code very easy to change—we just need to set the directory to look in and the target file extension.</w:t></w:r></w:p>'
Replace-ParagraphByAnchor "Because we were developing in IntelliJ" $notes2Xml
